{"js": "const replacements = [\n  [\"2025-02-28 Friday\", \"2025-03-01 Saturday\"],\n  [\"67\u00d737=\", \"24\u00d783=\"],\n  [\"25\u00d778=\", \"93\u00d781=\"],\n  [\"31\u00d745=\", \"40\u00d742=\"],\n  [\"23\u00d732=\", \"92\u00d748=\"],\n  [\"24\u00d771=\", \"97\u00d711=\"],\n  [\"97\u00d717=\", \"49\u00d773=\"],\n  [\"33\u00d736=\", \"36\u00d774=\"],\n  [\"32\u00d771=\", \"76\u00d767=\"],\n  [\"21\u00d740=\", \"90\u00d785=\"],\n  [\"51\u00d786=\", \"28\u00d764=\"],\n  [\"88\u00d753=\", \"69\u00d794=\"],\n  [\"34\u00d720=\", \"61\u00d747=\"],\n  [\"15\u00d767=\", \"26\u00d789=\"],\n  [\"67\u00d749=\", \"33\u00d731=\"],\n  [\"47\u00d716=\", \"27\u00d797=\"],\n  [\"79\u00d798=\", \"95\u00d746=\"],\n  [\"82\u00d719=\", \"54\u00d714=\"],\n  [\"39\u00d762=\", \"56\u00d720=\"],\n  [\"17\u00d731=\", \"48\u00d784=\"],\n  [\"64\u00d743=\", \"15\u00d788=\"],\n  [\"69\u00d773=\", \"69\u00d737=\"],\n  [\"35\u00d714=\", \"66\u00d787=\"],\n  [\"20\u00d753=\", \"23\u00d738=\"],\n  [\"39\u00d734=\", \"38\u00d745=\"],\n  [\"20\u00d717=\", \"28\u00d799=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-02-28 Friday\", \"2025-03-01 Saturday\"),\n    @(\"67\u00d737=\", \"24\u00d783=\"),\n    @(\"25\u00d778=\", \"93\u00d781=\"),\n    @(\"31\u00d745=\", \"40\u00d742=\"),\n    @(\"23\u00d732=\", \"92\u00d748=\"),\n    @(\"24\u00d771=\", \"97\u00d711=\"),\n    @(\"97\u00d717=\", \"49\u00d773=\"),\n    @(\"33\u00d736=\", \"36\u00d774=\"),\n    @(\"32\u00d771=\", \"76\u00d767=\"),\n    @(\"21\u00d740=\", \"90\u00d785=\"),\n    @(\"51\u00d786=\", \"28\u00d764=\"),\n    @(\"88\u00d753=\", \"69\u00d794=\"),\n    @(\"34\u00d720=\", \"61\u00d747=\"),\n    @(\"15\u00d767=\", \"26\u00d789=\"),\n    @(\"67\u00d749=\", \"33\u00d731=\"),\n    @(\"47\u00d716=\", \"27\u00d797=\"),\n    @(\"79\u00d798=\", \"95\u00d746=\"),\n    @(\"82\u00d719=\", \"54\u00d714=\"),\n    @(\"39\u00d762=\", \"56\u00d720=\"),\n    @(\"17\u00d731=\", \"48\u00d784=\"),\n    @(\"64\u00d743=\", \"15\u00d788=\"),\n    @(\"69\u00d773=\", \"69\u00d737=\"),\n    @(\"35\u00d714=\", \"66\u00d787=\"),\n    @(\"20\u00d753=\", \"23\u00d738=\"),\n    @(\"39\u00d734=\", \"38\u00d745=\"),\n    @(\"20\u00d717=\", \"28\u00d799=\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair[1]\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n}\n"}
